$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("B13").Value = 41986
$ws.Range("C13").Value = 0.64930555555555558
$ws.Range("D13").Value = 0.76388888888888884
$ws.Range("E13").Value = "sever"

# Row 14
$ws.Range("C14").Value = 0.79513888888888884
$ws.Range("D14").Value = 0.80902777777777779
$ws.Range("E14").Value = "sever - longtable"

# Row 15
$ws.Range("C15").Value = 0.89236111111111116
$ws.Range("D15").Value = 0.90625
$ws.Range("E15").Value = "schneider"

# Row 16
$ws.Range("C16").Value = 0.98958333333333337
$ws.Range("D16").Value = 0.052083333333333336
$ws.Range("E16").Value = "schneider"

# Row 17
$ws.Range("B17").Value = 41987
$ws.Range("C17").Value = 0.57291666666666663
$ws.Range("D17").Value = 0.625
$ws.Range("E17").Value = "schneider"

# Row 18
$ws.Range("C18").Value = 0.64930555555555558
$ws.Range("D18").Value = 0.69791666666666663
$ws.Range("E18").Value = "schneider tabs"

# Row 19
$ws.Range("C19").Value = 0.73263888888888884
$ws.Range("D19").Value = 0.78125
$ws.Range("E19").Value = "schneider images + text"

# Row 20
$ws.Range("C20").Value = 0.84027777777777779
$ws.Range("D20").Value = 0.86111111111111116
$ws.Range("E20").Value = "tiraz + submission"

# Row 21
$ws.Range("C21").Value = 0.94791666666666663
$ws.Range("D21").Value = 0.99305555555555547
$ws.Range("E21").Value = "Funta + Mako lit"

# Row 22
$ws.Range("C22").Value = 0.052083333333333336
$ws.Range("D22").Value = 0.072916666666666671
$ws.Range("E22").Value = "Sever lit"

# Extend the existing shared formula (D-C) from F5:F12 down to F5:F22,
# reproducing it as a single shared-formula group as Excel's fill-down would.
$ws.Range("F5:F22").Formula = "=D5-C5"

# Rows 12 and 16 keep literal values instead of the shared formula.
$ws.Range("F12").Value = 0.10069444444444443
$ws.Range("F16").Value = 0.0625

# Restore the selection to match the authored state.
$ws.Range("E23").Select()
